$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4
$ws.Range("A4").Value = "248F-19 "
$ws.Range("B4").Value = 20201027006
$ws.Range("C4").Value = 816.3610747358747
$ws.Range("D4").Value = 419.508702076143
$ws.Range("E4").Value = 2.533678884040108
$ws.Range("F4").Value = 0.7066148770379989

# Update row 5
$ws.Range("A5").Value = "248F-19 "
$ws.Range("B5").Value = 20201027007
$ws.Range("C5").Value = 171.8897684182772
$ws.Range("D5").Value = 0.06674757216893125
$ws.Range("E5").Value = 40.37504816535063
$ws.Range("F5").Value = 0.6831684792094123

# Add new row 6
$ws.Range("A6").Value = "248F-19 "
$ws.Range("B6").Value = 20201027008
$ws.Range("C6").Value = 172.8382378878981
$ws.Range("D6").Value = 0.01968585366899952
$ws.Range("E6").Value = 31.56235731912708
$ws.Range("F6").Value = 0.1959617602588581
